$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common (unchanged) columns for every data row in this dataset
$colA = 1
$colB = 'Agrícola del Norte S.A. de Arica'
$colC = 'Arica y Parinacota'
$colE = 15
$colF = 'Fruta'
$colG = 100104
$colH = 'Frutos de pepita'
$colI = 100104005
$colJ = 'Pera'

# Row 44
$ws.Cells.Item(44, 1).Value = $colA
$ws.Cells.Item(44, 2).Value = $colB
$ws.Cells.Item(44, 3).Value = $colC
$ws.Cells.Item(44, 4).Value = 45086
$ws.Cells.Item(44, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(44, 5).Value = $colE
$ws.Cells.Item(44, 6).Value = $colF
$ws.Cells.Item(44, 7).Value = $colG
$ws.Cells.Item(44, 8).Value = $colH
$ws.Cells.Item(44, 9).Value = $colI
$ws.Cells.Item(44, 10).Value = $colJ
$ws.Cells.Item(44, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(44, 12).Value = 'Segunda'
$ws.Cells.Item(44, 13).Value = 300
$ws.Cells.Item(44, 14).Value = 19000
$ws.Cells.Item(44, 15).Value = 20000
$ws.Cells.Item(44, 16).Value = 19500
$ws.Cells.Item(44, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(44, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(44, 19).Value = 1083
$ws.Cells.Item(44, 20).Value = 18

# Row 45
$ws.Cells.Item(45, 1).Value = $colA
$ws.Cells.Item(45, 2).Value = $colB
$ws.Cells.Item(45, 3).Value = $colC
$ws.Cells.Item(45, 4).Value = 45086
$ws.Cells.Item(45, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(45, 5).Value = $colE
$ws.Cells.Item(45, 6).Value = $colF
$ws.Cells.Item(45, 7).Value = $colG
$ws.Cells.Item(45, 8).Value = $colH
$ws.Cells.Item(45, 9).Value = $colI
$ws.Cells.Item(45, 10).Value = $colJ
$ws.Cells.Item(45, 11).Value = 'Winter Nelis'
$ws.Cells.Item(45, 12).Value = 'Segunda'
$ws.Cells.Item(45, 13).Value = 300
$ws.Cells.Item(45, 14).Value = 19000
$ws.Cells.Item(45, 15).Value = 20000
$ws.Cells.Item(45, 16).Value = 19500
$ws.Cells.Item(45, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(45, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(45, 19).Value = 1083
$ws.Cells.Item(45, 20).Value = 18

# Row 46
$ws.Cells.Item(46, 1).Value = $colA
$ws.Cells.Item(46, 2).Value = $colB
$ws.Cells.Item(46, 3).Value = $colC
$ws.Cells.Item(46, 4).Value = 44314
$ws.Cells.Item(46, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(46, 5).Value = $colE
$ws.Cells.Item(46, 6).Value = $colF
$ws.Cells.Item(46, 7).Value = $colG
$ws.Cells.Item(46, 8).Value = $colH
$ws.Cells.Item(46, 9).Value = $colI
$ws.Cells.Item(46, 10).Value = $colJ
$ws.Cells.Item(46, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(46, 12).Value = 'Segunda'
$ws.Cells.Item(46, 13).Value = 250
$ws.Cells.Item(46, 14).Value = 17000
$ws.Cells.Item(46, 15).Value = 18000
$ws.Cells.Item(46, 16).Value = 17500
$ws.Cells.Item(46, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(46, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(46, 19).Value = 972
$ws.Cells.Item(46, 20).Value = 18

# Row 47
$ws.Cells.Item(47, 1).Value = $colA
$ws.Cells.Item(47, 2).Value = $colB
$ws.Cells.Item(47, 3).Value = $colC
$ws.Cells.Item(47, 4).Value = 45072
$ws.Cells.Item(47, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(47, 5).Value = $colE
$ws.Cells.Item(47, 6).Value = $colF
$ws.Cells.Item(47, 7).Value = $colG
$ws.Cells.Item(47, 8).Value = $colH
$ws.Cells.Item(47, 9).Value = $colI
$ws.Cells.Item(47, 10).Value = $colJ
$ws.Cells.Item(47, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(47, 12).Value = 'Primera'
$ws.Cells.Item(47, 13).Value = 450
$ws.Cells.Item(47, 14).Value = 19000
$ws.Cells.Item(47, 15).Value = 20000
$ws.Cells.Item(47, 16).Value = 19556
$ws.Cells.Item(47, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(47, 18).Value = 'Región del Maule'
$ws.Cells.Item(47, 19).Value = 1086
$ws.Cells.Item(47, 20).Value = 18

# Row 48
$ws.Cells.Item(48, 1).Value = $colA
$ws.Cells.Item(48, 2).Value = $colB
$ws.Cells.Item(48, 3).Value = $colC
$ws.Cells.Item(48, 4).Value = 45072
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value = $colE
$ws.Cells.Item(48, 6).Value = $colF
$ws.Cells.Item(48, 7).Value = $colG
$ws.Cells.Item(48, 8).Value = $colH
$ws.Cells.Item(48, 9).Value = $colI
$ws.Cells.Item(48, 10).Value = $colJ
$ws.Cells.Item(48, 11).Value = 'Winter Nelis'
$ws.Cells.Item(48, 12).Value = 'Segunda'
$ws.Cells.Item(48, 13).Value = 350
$ws.Cells.Item(48, 14).Value = 19000
$ws.Cells.Item(48, 15).Value = 20000
$ws.Cells.Item(48, 16).Value = 19429
$ws.Cells.Item(48, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(48, 18).Value = 'Región del Maule'
$ws.Cells.Item(48, 19).Value = 1079
$ws.Cells.Item(48, 20).Value = 18

# Row 49
$ws.Cells.Item(49, 1).Value = $colA
$ws.Cells.Item(49, 2).Value = $colB
$ws.Cells.Item(49, 3).Value = $colC
$ws.Cells.Item(49, 4).Value = 44791
$ws.Cells.Item(49, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(49, 5).Value = $colE
$ws.Cells.Item(49, 6).Value = $colF
$ws.Cells.Item(49, 7).Value = $colG
$ws.Cells.Item(49, 8).Value = $colH
$ws.Cells.Item(49, 9).Value = $colI
$ws.Cells.Item(49, 10).Value = $colJ
$ws.Cells.Item(49, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(49, 12).Value = 'Primera'
$ws.Cells.Item(49, 13).Value = 250
$ws.Cells.Item(49, 14).Value = 19000
$ws.Cells.Item(49, 15).Value = 20000
$ws.Cells.Item(49, 16).Value = 19500
$ws.Cells.Item(49, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(49, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(49, 19).Value = 1083
$ws.Cells.Item(49, 20).Value = 18

# Row 50
$ws.Cells.Item(50, 1).Value = $colA
$ws.Cells.Item(50, 2).Value = $colB
$ws.Cells.Item(50, 3).Value = $colC
$ws.Cells.Item(50, 4).Value = 44791
$ws.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 5).Value = $colE
$ws.Cells.Item(50, 6).Value = $colF
$ws.Cells.Item(50, 7).Value = $colG
$ws.Cells.Item(50, 8).Value = $colH
$ws.Cells.Item(50, 9).Value = $colI
$ws.Cells.Item(50, 10).Value = $colJ
$ws.Cells.Item(50, 11).Value = 'Winter Nelis'
$ws.Cells.Item(50, 12).Value = 'Primera'
$ws.Cells.Item(50, 13).Value = 270
$ws.Cells.Item(50, 14).Value = 19000
$ws.Cells.Item(50, 15).Value = 20000
$ws.Cells.Item(50, 16).Value = 19500
$ws.Cells.Item(50, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(50, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(50, 19).Value = 1083
$ws.Cells.Item(50, 20).Value = 18

# Row 51
$ws.Cells.Item(51, 1).Value = $colA
$ws.Cells.Item(51, 2).Value = $colB
$ws.Cells.Item(51, 3).Value = $colC
$ws.Cells.Item(51, 4).Value = 44355
$ws.Cells.Item(51, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(51, 5).Value = $colE
$ws.Cells.Item(51, 6).Value = $colF
$ws.Cells.Item(51, 7).Value = $colG
$ws.Cells.Item(51, 8).Value = $colH
$ws.Cells.Item(51, 9).Value = $colI
$ws.Cells.Item(51, 10).Value = $colJ
$ws.Cells.Item(51, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(51, 12).Value = 'Segunda'
$ws.Cells.Item(51, 13).Value = 200
$ws.Cells.Item(51, 14).Value = 17000
$ws.Cells.Item(51, 15).Value = 18000
$ws.Cells.Item(51, 16).Value = 17500
$ws.Cells.Item(51, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(51, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(51, 19).Value = 972
$ws.Cells.Item(51, 20).Value = 18

# Row 52
$ws.Cells.Item(52, 1).Value = $colA
$ws.Cells.Item(52, 2).Value = $colB
$ws.Cells.Item(52, 3).Value = $colC
$ws.Cells.Item(52, 4).Value = 44355
$ws.Cells.Item(52, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(52, 5).Value = $colE
$ws.Cells.Item(52, 6).Value = $colF
$ws.Cells.Item(52, 7).Value = $colG
$ws.Cells.Item(52, 8).Value = $colH
$ws.Cells.Item(52, 9).Value = $colI
$ws.Cells.Item(52, 10).Value = $colJ
$ws.Cells.Item(52, 11).Value = 'Winter Nelis'
$ws.Cells.Item(52, 12).Value = 'Segunda'
$ws.Cells.Item(52, 13).Value = 250
$ws.Cells.Item(52, 14).Value = 17000
$ws.Cells.Item(52, 15).Value = 18000
$ws.Cells.Item(52, 16).Value = 17500
$ws.Cells.Item(52, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(52, 18).Value = 'Región Metropolitana'
$ws.Cells.Item(52, 19).Value = 972
$ws.Cells.Item(52, 20).Value = 18

# Row 53
$ws.Cells.Item(53, 1).Value = $colA
$ws.Cells.Item(53, 2).Value = $colB
$ws.Cells.Item(53, 3).Value = $colC
$ws.Cells.Item(53, 4).Value = 44336
$ws.Cells.Item(53, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(53, 5).Value = $colE
$ws.Cells.Item(53, 6).Value = $colF
$ws.Cells.Item(53, 7).Value = $colG
$ws.Cells.Item(53, 8).Value = $colH
$ws.Cells.Item(53, 9).Value = $colI
$ws.Cells.Item(53, 10).Value = $colJ
$ws.Cells.Item(53, 11).Value = 'Winter Nelis'
$ws.Cells.Item(53, 12).Value = 'Segunda'
$ws.Cells.Item(53, 13).Value = 250
$ws.Cells.Item(53, 14).Value = 21000
$ws.Cells.Item(53, 15).Value = 22000
$ws.Cells.Item(53, 16).Value = 21500
$ws.Cells.Item(53, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(53, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(53, 19).Value = 1194
$ws.Cells.Item(53, 20).Value = 18

# Row 54
$ws.Cells.Item(54, 1).Value = $colA
$ws.Cells.Item(54, 2).Value = $colB
$ws.Cells.Item(54, 3).Value = $colC
$ws.Cells.Item(54, 4).Value = 44994
$ws.Cells.Item(54, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(54, 5).Value = $colE
$ws.Cells.Item(54, 6).Value = $colF
$ws.Cells.Item(54, 7).Value = $colG
$ws.Cells.Item(54, 8).Value = $colH
$ws.Cells.Item(54, 9).Value = $colI
$ws.Cells.Item(54, 10).Value = $colJ
$ws.Cells.Item(54, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(54, 12).Value = 'Segunda'
$ws.Cells.Item(54, 13).Value = 150
$ws.Cells.Item(54, 14).Value = 19000
$ws.Cells.Item(54, 15).Value = 20000
$ws.Cells.Item(54, 16).Value = 19333
$ws.Cells.Item(54, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(54, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(54, 19).Value = 1074
$ws.Cells.Item(54, 20).Value = 18

# Row 55
$ws.Cells.Item(55, 1).Value = $colA
$ws.Cells.Item(55, 2).Value = $colB
$ws.Cells.Item(55, 3).Value = $colC
$ws.Cells.Item(55, 4).Value = 44497
$ws.Cells.Item(55, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(55, 5).Value = $colE
$ws.Cells.Item(55, 6).Value = $colF
$ws.Cells.Item(55, 7).Value = $colG
$ws.Cells.Item(55, 8).Value = $colH
$ws.Cells.Item(55, 9).Value = $colI
$ws.Cells.Item(55, 10).Value = $colJ
$ws.Cells.Item(55, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(55, 12).Value = 'Segunda'
$ws.Cells.Item(55, 13).Value = 300
$ws.Cells.Item(55, 14).Value = 17000
$ws.Cells.Item(55, 15).Value = 18000
$ws.Cells.Item(55, 16).Value = 17500
$ws.Cells.Item(55, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(55, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(55, 19).Value = 972
$ws.Cells.Item(55, 20).Value = 18

# Row 56
$ws.Cells.Item(56, 1).Value = $colA
$ws.Cells.Item(56, 2).Value = $colB
$ws.Cells.Item(56, 3).Value = $colC
$ws.Cells.Item(56, 4).Value = 44497
$ws.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 5).Value = $colE
$ws.Cells.Item(56, 6).Value = $colF
$ws.Cells.Item(56, 7).Value = $colG
$ws.Cells.Item(56, 8).Value = $colH
$ws.Cells.Item(56, 9).Value = $colI
$ws.Cells.Item(56, 10).Value = $colJ
$ws.Cells.Item(56, 11).Value = 'Winter Nelis'
$ws.Cells.Item(56, 12).Value = 'Segunda'
$ws.Cells.Item(56, 13).Value = 250
$ws.Cells.Item(56, 14).Value = 17000
$ws.Cells.Item(56, 15).Value = 18000
$ws.Cells.Item(56, 16).Value = 17500
$ws.Cells.Item(56, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(56, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(56, 19).Value = 972
$ws.Cells.Item(56, 20).Value = 18

# Row 57
$ws.Cells.Item(57, 1).Value = $colA
$ws.Cells.Item(57, 2).Value = $colB
$ws.Cells.Item(57, 3).Value = $colC
$ws.Cells.Item(57, 4).Value = 44883
$ws.Cells.Item(57, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(57, 5).Value = $colE
$ws.Cells.Item(57, 6).Value = $colF
$ws.Cells.Item(57, 7).Value = $colG
$ws.Cells.Item(57, 8).Value = $colH
$ws.Cells.Item(57, 9).Value = $colI
$ws.Cells.Item(57, 10).Value = $colJ
$ws.Cells.Item(57, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(57, 12).Value = 'Primera'
$ws.Cells.Item(57, 13).Value = 300
$ws.Cells.Item(57, 14).Value = 24000
$ws.Cells.Item(57, 15).Value = 25000
$ws.Cells.Item(57, 16).Value = 24500
$ws.Cells.Item(57, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(57, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(57, 19).Value = 1361
$ws.Cells.Item(57, 20).Value = 18

# Row 58
$ws.Cells.Item(58, 1).Value = $colA
$ws.Cells.Item(58, 2).Value = $colB
$ws.Cells.Item(58, 3).Value = $colC
$ws.Cells.Item(58, 4).Value = 44525
$ws.Cells.Item(58, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(58, 5).Value = $colE
$ws.Cells.Item(58, 6).Value = $colF
$ws.Cells.Item(58, 7).Value = $colG
$ws.Cells.Item(58, 8).Value = $colH
$ws.Cells.Item(58, 9).Value = $colI
$ws.Cells.Item(58, 10).Value = $colJ
$ws.Cells.Item(58, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(58, 12).Value = 'Segunda'
$ws.Cells.Item(58, 13).Value = 300
$ws.Cells.Item(58, 14).Value = 19000
$ws.Cells.Item(58, 15).Value = 20000
$ws.Cells.Item(58, 16).Value = 19500
$ws.Cells.Item(58, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(58, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(58, 19).Value = 1083
$ws.Cells.Item(58, 20).Value = 18

# Row 59
$ws.Cells.Item(59, 1).Value = $colA
$ws.Cells.Item(59, 2).Value = $colB
$ws.Cells.Item(59, 3).Value = $colC
$ws.Cells.Item(59, 4).Value = 44993
$ws.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(59, 5).Value = $colE
$ws.Cells.Item(59, 6).Value = $colF
$ws.Cells.Item(59, 7).Value = $colG
$ws.Cells.Item(59, 8).Value = $colH
$ws.Cells.Item(59, 9).Value = $colI
$ws.Cells.Item(59, 10).Value = $colJ
$ws.Cells.Item(59, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(59, 12).Value = 'Primera'
$ws.Cells.Item(59, 13).Value = 550
$ws.Cells.Item(59, 14).Value = 14500
$ws.Cells.Item(59, 15).Value = 15000
$ws.Cells.Item(59, 16).Value = 14727
$ws.Cells.Item(59, 17).Value = '$/caja 20 kilos granel'
$ws.Cells.Item(59, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(59, 19).Value = 736
$ws.Cells.Item(59, 20).Value = 20

# Row 60
$ws.Cells.Item(60, 1).Value = $colA
$ws.Cells.Item(60, 2).Value = $colB
$ws.Cells.Item(60, 3).Value = $colC
$ws.Cells.Item(60, 4).Value = 44421
$ws.Cells.Item(60, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(60, 5).Value = $colE
$ws.Cells.Item(60, 6).Value = $colF
$ws.Cells.Item(60, 7).Value = $colG
$ws.Cells.Item(60, 8).Value = $colH
$ws.Cells.Item(60, 9).Value = $colI
$ws.Cells.Item(60, 10).Value = $colJ
$ws.Cells.Item(60, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(60, 12).Value = 'Segunda'
$ws.Cells.Item(60, 13).Value = 270
$ws.Cells.Item(60, 14).Value = 16000
$ws.Cells.Item(60, 15).Value = 17000
$ws.Cells.Item(60, 16).Value = 16500
$ws.Cells.Item(60, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(60, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(60, 19).Value = 917
$ws.Cells.Item(60, 20).Value = 18

# Row 61
$ws.Cells.Item(61, 1).Value = $colA
$ws.Cells.Item(61, 2).Value = $colB
$ws.Cells.Item(61, 3).Value = $colC
$ws.Cells.Item(61, 4).Value = 44421
$ws.Cells.Item(61, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(61, 5).Value = $colE
$ws.Cells.Item(61, 6).Value = $colF
$ws.Cells.Item(61, 7).Value = $colG
$ws.Cells.Item(61, 8).Value = $colH
$ws.Cells.Item(61, 9).Value = $colI
$ws.Cells.Item(61, 10).Value = $colJ
$ws.Cells.Item(61, 11).Value = 'Winter Nelis'
$ws.Cells.Item(61, 12).Value = 'Segunda'
$ws.Cells.Item(61, 13).Value = 250
$ws.Cells.Item(61, 14).Value = 16000
$ws.Cells.Item(61, 15).Value = 17000
$ws.Cells.Item(61, 16).Value = 16500
$ws.Cells.Item(61, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(61, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(61, 19).Value = 917
$ws.Cells.Item(61, 20).Value = 18

# Row 62
$ws.Cells.Item(62, 1).Value = $colA
$ws.Cells.Item(62, 2).Value = $colB
$ws.Cells.Item(62, 3).Value = $colC
$ws.Cells.Item(62, 4).Value = 45021
$ws.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 5).Value = $colE
$ws.Cells.Item(62, 6).Value = $colF
$ws.Cells.Item(62, 7).Value = $colG
$ws.Cells.Item(62, 8).Value = $colH
$ws.Cells.Item(62, 9).Value = $colI
$ws.Cells.Item(62, 10).Value = $colJ
$ws.Cells.Item(62, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(62, 12).Value = 'Segunda'
$ws.Cells.Item(62, 13).Value = 300
$ws.Cells.Item(62, 14).Value = 19000
$ws.Cells.Item(62, 15).Value = 20000
$ws.Cells.Item(62, 16).Value = 19500
$ws.Cells.Item(62, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(62, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(62, 19).Value = 1083
$ws.Cells.Item(62, 20).Value = 18

# Row 63
$ws.Cells.Item(63, 1).Value = $colA
$ws.Cells.Item(63, 2).Value = $colB
$ws.Cells.Item(63, 3).Value = $colC
$ws.Cells.Item(63, 4).Value = 45021
$ws.Cells.Item(63, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 5).Value = $colE
$ws.Cells.Item(63, 6).Value = $colF
$ws.Cells.Item(63, 7).Value = $colG
$ws.Cells.Item(63, 8).Value = $colH
$ws.Cells.Item(63, 9).Value = $colI
$ws.Cells.Item(63, 10).Value = $colJ
$ws.Cells.Item(63, 11).Value = 'Winter Nelis'
$ws.Cells.Item(63, 12).Value = 'Segunda'
$ws.Cells.Item(63, 13).Value = 300
$ws.Cells.Item(63, 14).Value = 19000
$ws.Cells.Item(63, 15).Value = 20000
$ws.Cells.Item(63, 16).Value = 19500
$ws.Cells.Item(63, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(63, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(63, 19).Value = 1083
$ws.Cells.Item(63, 20).Value = 18

# Row 64
$ws.Cells.Item(64, 1).Value = $colA
$ws.Cells.Item(64, 2).Value = $colB
$ws.Cells.Item(64, 3).Value = $colC
$ws.Cells.Item(64, 4).Value = 44371
$ws.Cells.Item(64, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(64, 5).Value = $colE
$ws.Cells.Item(64, 6).Value = $colF
$ws.Cells.Item(64, 7).Value = $colG
$ws.Cells.Item(64, 8).Value = $colH
$ws.Cells.Item(64, 9).Value = $colI
$ws.Cells.Item(64, 10).Value = $colJ
$ws.Cells.Item(64, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(64, 12).Value = 'Calibre 90'
$ws.Cells.Item(64, 13).Value = 140
$ws.Cells.Item(64, 14).Value = 17000
$ws.Cells.Item(64, 15).Value = 18000
$ws.Cells.Item(64, 16).Value = 17429
$ws.Cells.Item(64, 17).Value = '$/caja 18 kilos embalada'
$ws.Cells.Item(64, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(64, 19).Value = 968
$ws.Cells.Item(64, 20).Value = 18

# Row 65
$ws.Cells.Item(65, 1).Value = $colA
$ws.Cells.Item(65, 2).Value = $colB
$ws.Cells.Item(65, 3).Value = $colC
$ws.Cells.Item(65, 4).Value = 44371
$ws.Cells.Item(65, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(65, 5).Value = $colE
$ws.Cells.Item(65, 6).Value = $colF
$ws.Cells.Item(65, 7).Value = $colG
$ws.Cells.Item(65, 8).Value = $colH
$ws.Cells.Item(65, 9).Value = $colI
$ws.Cells.Item(65, 10).Value = $colJ
$ws.Cells.Item(65, 11).Value = 'Winter Nelis'
$ws.Cells.Item(65, 12).Value = 'Calibre 80'
$ws.Cells.Item(65, 13).Value = 120
$ws.Cells.Item(65, 14).Value = 17000
$ws.Cells.Item(65, 15).Value = 18000
$ws.Cells.Item(65, 16).Value = 17500
$ws.Cells.Item(65, 17).Value = '$/caja 18 kilos embalada'
$ws.Cells.Item(65, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(65, 19).Value = 972
$ws.Cells.Item(65, 20).Value = 18

# Row 66
$ws.Cells.Item(66, 1).Value = $colA
$ws.Cells.Item(66, 2).Value = $colB
$ws.Cells.Item(66, 3).Value = $colC
$ws.Cells.Item(66, 4).Value = 44981
$ws.Cells.Item(66, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(66, 5).Value = $colE
$ws.Cells.Item(66, 6).Value = $colF
$ws.Cells.Item(66, 7).Value = $colG
$ws.Cells.Item(66, 8).Value = $colH
$ws.Cells.Item(66, 9).Value = $colI
$ws.Cells.Item(66, 10).Value = $colJ
$ws.Cells.Item(66, 11).Value = 'Winter Nelis'
$ws.Cells.Item(66, 12).Value = 'Segunda'
$ws.Cells.Item(66, 13).Value = 300
$ws.Cells.Item(66, 14).Value = 25000
$ws.Cells.Item(66, 15).Value = 26000
$ws.Cells.Item(66, 16).Value = 25500
$ws.Cells.Item(66, 17).Value = '$/caja 20 kilos empedrada'
$ws.Cells.Item(66, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(66, 19).Value = 1275
$ws.Cells.Item(66, 20).Value = 20

Write-Output "Rows 44-66 updated"